$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.23%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.53%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.201"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.46%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07683"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.42%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.295"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.64%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.706"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.40%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9405"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.44%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.84%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1832"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.49%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09157"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.64%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04244"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.20%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1053"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.94%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001283"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.06%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005891"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.03%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.350"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.18%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.532"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "8.72%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1344"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2722"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.52%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04018"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.10%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001267"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.34%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004240"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.18%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001271"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.07%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02536"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.17%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.85%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007840"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.53%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1314"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.82%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006670"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.71%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001942"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.40%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008094"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.24%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3092"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.45%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006773"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.15%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2210"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "174.52%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "3.42%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
